$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "Choice A"
$ws.Range("E1").Value = "Choice B"
$ws.Range("F1").Value = "Choice C"
$ws.Range("G1").Value = "Choice D"

$ws.Range("F1").Select() | Out-Null
